$d = $word.ActiveDocument

# --- Fix typo: "Virtual Netorking" -> "Virtual Networking" -----------------
$null = $d.Content.Find.Execute("Virtual Netorking", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Virtual Networking", 2)

# --- Locate the last bullet ("Virtual Networking") to append new bullets ---
$lastBullet = $d.Paragraphs.Item($d.Paragraphs.Count).Range

# New top-level bullet: Lauren Tureaud- Vandy DA consultant (Operations side)
$lastBullet.InsertParagraphAfter()
$lauren = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lauren.ListFormat.ListLevelNumber = 1
$lauren.Text = "Lauren Tureaud- Vandy DA consultant (Operations side)"

# Sub-bullet: Azure
$lauren.InsertParagraphAfter()
$azure = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$azure.ListFormat.ListLevelNumber = 2
$azure.Text = "Azure"

# Sub-bullet: SQL, Tableau
$azure.InsertParagraphAfter()
$sql = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$sql.ListFormat.ListLevelNumber = 2
$sql.Text = "SQL, Tableau "

# Sub-bullet: Microsoft teams
$sql.InsertParagraphAfter()
$teams = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$teams.ListFormat.ListLevelNumber = 2
$teams.Text = "Microsoft teams "

# --- Trailing blank paragraph at the very end of the document --------------
$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")
